$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style_D2 = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.844.43'
$ws.Range('D2').Style = $style_D2
$ws.Range('E2').Value = '  +0.27%  '
$style_D3 = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.085.01'
$ws.Range('D3').Style = $style_D3
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('E4').Value = '  +0.13%  '
$style_D5 = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.42'
$ws.Range('D5').Style = $style_D5
$ws.Range('E5').Value = '  +0.73%  '
$style_D6 = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.628'
$ws.Range('D6').Style = $style_D6
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('E7').Value = '  +3.43%  '
$ws.Range('E8').Value = '  +0.02%  '
$style_D9 = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.395'
$ws.Range('D9').Style = $style_D9
$ws.Range('E9').Value = '  +2.14%  '
$ws.Range('E11').Value = '  +3.91%  '
$style_D12 = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.391.69'
$ws.Range('D12').Style = $style_D12
$ws.Range('E12').Value = '  +0.25%  '
$style_D13 = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.79'
$ws.Range('D13').Style = $style_D13
$ws.Range('E13').Value = '  +1.96%  '
$style_D14 = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.35'
$ws.Range('D14').Style = $style_D14
$ws.Range('E14').Value = '  +2.02%  '
$ws.Range('E15').Value = '  +2.62%  '
$ws.Range('E16').Value = '  +1.68%  '
$style_D17 = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.137.26'
$ws.Range('D17').Style = $style_D17
$ws.Range('E17').Value = '  +2.67%  '
$style_D18 = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.819.76'
$ws.Range('D18').Style = $style_D18
$ws.Range('E18').Value = '  +0.40%  '
$style_D19 = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.17'
$ws.Range('D19').Style = $style_D19
$ws.Range('E19').Value = '  +0.12%  '
$style_D20 = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.81'
$ws.Range('D20').Style = $style_D20
$ws.Range('E20').Value = '  +1.52%  '
$style_D21 = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0848'
$ws.Range('D21').Style = $style_D21
$style_D22 = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.61'
$ws.Range('D22').Style = $style_D22
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('E25').Value = '  +1.68%  '
$style_D26 = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.28'
$ws.Range('D26').Style = $style_D26
$ws.Range('E26').Value = '  +0.65%  '
$style_D27 = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.49'
$ws.Range('D27').Style = $style_D27
$ws.Range('E27').Value = '  +6.49%  '
$style_D28 = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.137'
$ws.Range('D28').Style = $style_D28
$ws.Range('E28').Value = '  -1.47%  '
$ws.Range('E29').Value = '  -0.80%  '
$style_D30 = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.57'
$ws.Range('D30').Style = $style_D30
$ws.Range('E30').Value = '  +0.93%  '
$ws.Range('E31').Value = '  +2.54%  '
$style_D32 = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.78'
$ws.Range('D32').Style = $style_D32
$ws.Range('E32').Value = '  +2.99%  '
$style_D33 = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0636'
$ws.Range('D33').Style = $style_D33
$ws.Range('E33').Value = '  +1.93%  '
$style_D34 = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.72'
$ws.Range('D34').Style = $style_D34
$ws.Range('E34').Value = '  +2.84%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('E36').Value = '  +2.22%  '
$ws.Range('E37').Value = '  -0.27%  '
$style_D38 = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').Style = $style_D38
$ws.Range('E38').Value = '  -0.07%  '
$style_D39 = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.44'
$ws.Range('D39').Style = $style_D39
$ws.Range('E39').Value = '  +0.12%  '
$style_D40 = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0982'
$ws.Range('D40').Style = $style_D40
$ws.Range('E40').Value = '  -0.47%  '
$style_D41 = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.51'
$ws.Range('D41').Style = $style_D41
$ws.Range('E41').Value = '  +0.52%  '
$style_D42 = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0222'
$ws.Range('D42').Style = $style_D42
$ws.Range('E42').Value = '  +3.82%  '
$style_D43 = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.14'
$ws.Range('D43').Style = $style_D43
$ws.Range('E43').Value = '  +9.89%  '
$ws.Range('E44').Value = '  -0.90%  '
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('E46').Value = '  -0.67%  '
$style_D47 = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.20'
$ws.Range('D47').Style = $style_D47
$ws.Range('E47').Value = '  -3.97%  '
$ws.Range('E48').Value = '  +1.64%  '
$style_D49 = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.40'
$ws.Range('D49').Style = $style_D49
$ws.Range('E49').Value = '  -0.49%  '
$ws.Range('E50').Value = '  -0.87%  '
$style_D51 = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.276.89'
$ws.Range('D51').Style = $style_D51
$ws.Range('E51').Value = '  +0.23%  '
